# Issue #19 select/unselect all in directory bar
#
# - Adds a new follow-up issue (#20) about returning from the images
#   browser to the playlist via the navbar.
# - Marks issue #19 ("select/unselect all") as DONE and renames it to
#   "select/unselect all in directory bar".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# New row 21 = Issue #20 "images to playlist return nav"
# (write the Description (H) before the Name (E) so new shared strings
# land in the same order as the source workbook).
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 1
$ws.Range("H21").Value = "add return from images to playlist on navbar"
$ws.Range("E21").Value = "images to playlist return nav"
$ws.Rows.Item(21).RowHeight = 29

# Issue #19 (row 20): close it out and give it its final name.
$ws.Range("C20").Value = "DONE"
$ws.Range("E20").Value = "select/unselect all in directory bar"

# Leave the selection where the editor ended up after adding the row.
$ws.Range("E25").Select()
